# Update "想去人数" (people interested) counts that changed between scrapes.
$wb = $excel.ActiveWorkbook

# Sheet "展览" (Exhibitions)
$wsExhibit = $wb.Worksheets.Item("展览")
$wsExhibit.Range("F2").Value = 384
$wsExhibit.Range("F3").Value = 2133
$wsExhibit.Range("F4").Value = 110

# Sheet "全部类型" (All types) mirrors the same rows from "展览"
$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F2").Value = 384
$wsAll.Range("F7").Value = 2133
$wsAll.Range("F8").Value = 110
